# The underlying edit captured by the diff is a simple re-save of the
# workbook in which the only deliberate, user-driven change is where the
# cursor/selection ended up on Sheet1: it moves from A2 to C7.
# (The surrounding noise in the diff -- fileVersion/rupBuild, the
# x15ac:absPath username/machine path, revisionPtr GUIDs/coauth versions,
# the bookViews window geometry, and the tiny defaultRowHeight / row
# ht / x14ac:dyDescent / column-width jitter -- is all incidental
# metadata that Excel re-derives from the authoring machine, Office
# build, and screen/font metrics whenever a file is opened and saved
# again; it isn't something a script drives on purpose.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Select()
